# Append the new daily row (2025/10/09, 木, 20, 199) to the bottom of the
# sheet1 table, mirroring the existing rows (A:date text, B:weekday text,
# C/D: plain numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 85
$lastRow = $newRow - 1

# Column A holds a date-looking string ("2025/09/22" etc.) that is stored as
# literal text in the workbook, not a real date. Assigning that text straight
# to .Value would make Excel auto-detect it as a date (changing the cell's
# type/format). Force the cell to Text first, write the value, then restore
# the plain/default style (matching the unstyled cells above it) so the
# written cell looks identical to its neighbours.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/10/09"
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($lastRow, 1).Style

$ws.Cells.Item($newRow, 2).Value = "木"
$ws.Cells.Item($newRow, 3).Value = 20
$ws.Cells.Item($newRow, 4).Value = 199
